$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.972438577363846
$ws.Range("C2").Value = 0.2623415173089541
$ws.Range("D2").Value = 0.5882220620855207
$ws.Range("E2").Value = 0.2184013329535901
$ws.Range("G2").Value = 0.002473343836612459
$ws.Range("I2").Value = 0.740465206117733
$ws.Range("J2").Value = 0.09951159869620341
$ws.Range("M2").Value = 0.4923926705346986
$ws.Range("O2").Value = 4.066081606316061
$ws.Range("B3").Value = 0.8753874355857079
$ws.Range("C3").Value = 0.2346213734578271
$ws.Range("D3").Value = 0.585045878658164
$ws.Range("E3").Value = 0.2184604951824554
$ws.Range("G3").Value = 0.002476804055579263
$ws.Range("I3").Value = 0.750562122363128
$ws.Range("J3").Value = 0.1002439673492326
$ws.Range("M3").Value = 0.465046541074571
$ws.Range("O3").Value = 4.085892049871745
$ws.Range("B4").Value = 0.8158348109444375
$ws.Range("C4").Value = 0.2176024792912301
$ws.Range("D4").Value = 0.5833877914308232
$ws.Range("E4").Value = 0.2186047646559146
$ws.Range("G4").Value = 0.002479041684627488
$ws.Range("I4").Value = 0.7572949739284702
$ws.Range("J4").Value = 0.1007405330698994
$ws.Range("M4").Value = 0.4484092753337805
$ws.Range("O4").Value = 4.10093991790734
$ws.Range("B5").Value = 0.7915772034498048
$ws.Range("C5").Value = 0.2106677657038745
$ws.Range("D5").Value = 0.5827855856817621
$ws.Range("E5").Value = 0.2186906979776602
$ws.Range("G5").Value = 0.002479982050146767
$ws.Range("I5").Value = 0.7601725124401355
$ws.Range("J5").Value = 0.1009546732079603
$ws.Range("M5").Value = 0.441668304562377
$ws.Range("O5").Value = 4.107795947460971
$ws.Range("B6").Value = 0.7875499210352928
$ws.Range("C6").Value = 0.2095163076937183
$ws.Range("D6").Value = 0.5826900286028831
$ws.Range("E6").Value = 0.2187066064519385
$ws.Range("G6").Value = 0.002480139922012993
$ws.Range("I6").Value = 0.7606584020137603
$ws.Range("J6").Value = 0.1009909426932083
$ws.Range("M6").Value = 0.4405513264270553
$ws.Range("O6").Value = 4.108978068362887
$ws.Range("B7").Value = 0.815507619931509
$ws.Range("C7").Value = 0.2175089524218663
$ws.Range("D7").Value = 0.5833793723252683
$ws.Range("E7").Value = 0.2186058136855067
$ws.Range("G7").Value = 0.002479054251062564
$ws.Range("I7").Value = 0.7573332398225574
$ws.Range("J7").Value = 0.1007433733261394
$ws.Range("M7").Value = 0.4483182065613107
$ws.Range("O7").Value = 4.101029451475853
$ws.Range("B8").Value = 0.9389683286097465
$ws.Range("C8").Value = 0.2527834793545765
$ws.Range("D8").Value = 0.5870663015025457
$ws.Range("E8").Value = 0.2183993271818743
$ws.Range("G8").Value = 0.002474513512016838
$ws.Range("I8").Value = 0.7438357816711907
$ws.Range("J8").Value = 0.09975438472355691
$ws.Range("M8").Value = 0.482932053651318
$ws.Range("O8").Value = 4.072312818256052
$ws.Range("B9").Value = 1.181328301003816
$ws.Range("C9").Value = 0.3219592485391445
$ws.Range("D9").Value = 0.5966140459414646
$ws.Range("E9").Value = 0.218851209818645
$ws.Range("G9").Value = 0.00246650201217051
$ws.Range("I9").Value = 0.7216102846072197
$ws.Range("J9").Value = 0.09818726704788361
$ws.Range("M9").Value = 0.5520180603939622
$ws.Range("O9").Value = 4.0389486988756
$ws.Range("B10").Value = 1.359508651215492
$ws.Range("C10").Value = 0.3727782026186333
$ws.Range("D10").Value = 0.6050434095509161
$ws.Range("E10").Value = 0.2197063751221613
$ws.Range("G10").Value = 0.002461154604383175
$ws.Range("I10").Value = 0.7078834984591289
$ws.Range("J10").Value = 0.09726324313883339
$ws.Range("M10").Value = 0.6035065312700283
$ws.Range("O10").Value = 4.02852101812482
$ws.Range("B11").Value = 1.440586952492538
$ws.Range("C11").Value = 0.3958952109169331
$ws.Range("D11").Value = 0.6091858181078749
$ws.Range("E11").Value = 0.2202092254044565
$ws.Range("G11").Value = 0.00245883768083648
$ws.Range("I11").Value = 0.7022072654228069
$ws.Range("J11").Value = 0.0968923332169318
$ws.Range("M11").Value = 0.6270879265664604
$ws.Range("O11").Value = 4.026856362535824
$ws.Range("B12").Value = 1.471291614716904
$ws.Range("C12").Value = 0.4046487223262147
$ws.Range("D12").Value = 0.6107987180110115
$ws.Range("E12").Value = 0.2204160215108324
$ws.Range("G12").Value = 0.00245797685947372
$ws.Range("I12").Value = 0.70013980947585
$ws.Range("J12").Value = 0.09675899517979758
$ws.Range("M12").Value = 0.6360402711756308
$ws.Range("O12").Value = 4.026670355924523
$ws.Range("B13").Value = 1.464678744214041
$ws.Range("C13").Value = 0.4027635192024945
$ws.Range("D13").Value = 0.6104493828230773
$ws.Range("E13").Value = 0.2203707557266874
$ws.Range("G13").Value = 0.002458161518058708
$ws.Range("I13").Value = 0.7005814207184002
$ws.Range("J13").Value = 0.09678739521302049
$ws.Range("M13").Value = 0.634111223761451
$ws.Range("O13").Value = 4.026690626921408
$ws.Range("B14").Value = 1.44311300923124
$ws.Range("C14").Value = 0.3966153771555696
$ws.Range("D14").Value = 0.6093176253701529
$ws.Range("E14").Value = 0.22022591036108
$ws.Range("G14").Value = 0.002458766529564707
$ws.Range("I14").Value = 0.7020355291866309
$ws.Range("J14").Value = 0.0968812207189309
$ws.Range("M14").Value = 0.6278239903942193
$ws.Range("O14").Value = 4.026832144727848
$ws.Range("B15").Value = 1.429903606161304
$ws.Range("C15").Value = 0.392849402443801
$ws.Range("D15").Value = 0.608630154930637
$ws.Range("E15").Value = 0.2201393215088707
$ws.Range("G15").Value = 0.002459139268900377
$ws.Range("I15").Value = 0.7029369020082541
$ws.Range("J15").Value = 0.09693961873027135
$ws.Range("M15").Value = 0.6239758108299611
$ws.Range("O15").Value = 4.026976745484546
$ws.Range("B16").Value = 1.354210351425536
$ws.Range("C16").Value = 0.3712674114615311
$ws.Range("D16").Value = 0.6047788873554225
$ws.Range("E16").Value = 0.2196758039965765
$ws.Range("G16").Value = 0.002461308340735313
$ws.Range("I16").Value = 0.7082659128524291
$ws.Range("J16").Value = 0.09728847846826127
$ws.Range("M16").Value = 0.6019686071547028
$ws.Range("O16").Value = 4.028691901449889
$ws.Range("B17").Value = 1.307780021232304
$ws.Range("C17").Value = 0.3580271530907453
$ws.Range("D17").Value = 0.6024951004513355
$ws.Range("E17").Value = 0.2194206121136659
$ws.Range("G17").Value = 0.002462668555780089
$ws.Range("I17").Value = 0.7116808265507686
$ws.Range("J17").Value = 0.09751516005171723
$ws.Range("M17").Value = 0.5885084270271363
$ws.Range("O17").Value = 4.030533833824819
$ws.Range("B18").Value = 1.281076778199292
$ws.Range("C18").Value = 0.350411649137925
$ws.Range("D18").Value = 0.6012105032647526
$ws.Range("E18").Value = 0.219284546286076
$ws.Range("G18").Value = 0.002463461804783126
$ws.Range("I18").Value = 0.713698453009215
$ws.Range("J18").Value = 0.0976501936378078
$ws.Range("M18").Value = 0.5807814730888339
$ws.Range("O18").Value = 4.031882973184736
$ws.Range("B19").Value = 1.272035944950176
$ws.Range("C19").Value = 0.3478331692133168
$ws.Range("D19").Value = 0.600780537772863
$ws.Range("E19").Value = 0.2192403165759735
$ws.Range("G19").Value = 0.002463732258029649
$ws.Range("I19").Value = 0.7143907587340621
$ws.Range("J19").Value = 0.09769671245879508
$ws.Range("M19").Value = 0.5781678418586296
$ws.Range("O19").Value = 4.032389477226644
$ws.Range("B20").Value = 1.312722384973711
$ws.Range("C20").Value = 0.3594366092290215
$ws.Range("D20").Value = 0.6027352145598854
$ws.Range("E20").Value = 0.2194466688250429
$ws.Range("G20").Value = 0.00246252263201001
$ws.Range("I20").Value = 0.7113117681572483
$ws.Range("J20").Value = 0.09749054784551703
$ws.Range("M20").Value = 0.5899397373111839
$ws.Range("O20").Value = 4.030307761531674
$ws.Range("B21").Value = 1.449447342173073
$ws.Range("C21").Value = 0.398421248552097
$ws.Range("D21").Value = 0.6096488488363434
$ws.Range("E21").Value = 0.2202680104003782
$ws.Range("G21").Value = 0.002458588374469306
$ws.Range("I21").Value = 0.7016061940648228
$ws.Range("J21").Value = 0.09685346865365219
$ws.Range("M21").Value = 0.6296700926318977
$ws.Range("O21").Value = 4.026778504780737
$ws.Range("B22").Value = 1.538816492752176
$ws.Range("C22").Value = 0.4238975158365861
$ws.Range("D22").Value = 0.6144252845682558
$ws.Range("E22").Value = 0.2209002702138001
$ws.Range("G22").Value = 0.00245611352134905
$ws.Range("I22").Value = 0.6957412038219175
$ws.Range("J22").Value = 0.09647859111821333
$ws.Range("M22").Value = 0.6557676718742584
$ws.Range("O22").Value = 4.027062765824155
$ws.Range("B23").Value = 1.491117856588971
$ws.Range("C23").Value = 0.4103006693142675
$ws.Range("D23").Value = 0.6118524086633101
$ws.Range("E23").Value = 0.2205540831446164
$ws.Range("G23").Value = 0.002457425602320331
$ws.Range("I23").Value = 0.6988276046602877
$ws.Range("J23").Value = 0.09667487089391713
$ws.Range("M23").Value = 0.6418269643763637
$ws.Range("O23").Value = 4.026673455136034
$ws.Range("B24").Value = 1.310487972846602
$ws.Range("C24").Value = 0.3587994050035945
$ws.Range("D24").Value = 0.602626570563956
$ws.Range("E24").Value = 0.219434855421003
$ws.Range("G24").Value = 0.002462588568980635
$ws.Range("I24").Value = 0.7114784501864619
$ws.Range("J24").Value = 0.09750166034899621
$ws.Range("M24").Value = 0.5892926061440136
$ws.Range("O24").Value = 4.030409064991858
$ws.Range("B25").Value = 1.115740146066287
$ws.Range("C25").Value = 0.3032457402188697
$ws.Range("D25").Value = 0.5937828307081077
$ws.Range("E25").Value = 0.2184604951824554
$ws.Range("G25").Value = 0.002468574335831778
$ws.Range("I25").Value = 0.7271670140547819
$ws.Range("J25").Value = 0.0985713191511195
$ws.Range("M25").Value = 0.5520180603939622
$ws.Range("O25").Value = 4.0455079530538
